# Reset changes to callODESolver to use SUNDIALS instead of ode15s
#
# Targets the "N3mr" sheet (renamed to "N3mOpen"): adds REV/FWD effector
# columns (C1:E1 headers, C9 -> D9 shift) and switches the reaction A9 from
# an irreversible arrow to a reversible one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N3mr")

# Move the existing "A[c], K[c]" value out of C9 so column C can become the
# new "Negative Effectors FWD" column; it now belongs under the REV side (D9).
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "A[c], K[c]"

# Header row: keep B1 ("Positive Effectors") -> now "FWD", and add the three
# new effector columns.
$ws.Range("B1").Value = "Positive Effectors FWD"
$ws.Range("C1").Value = "Negative Effectors FWD"
$ws.Range("D1").Value = "Positive Effectors REV"
$ws.Range("E1").Value = "Negative Effectors REV"

# Switch the reaction text from irreversible to reversible.
$ws.Range("A9").Value = "M[c] + N[c] <==>A[c] + K[c]"

# Update selection / active cell to match the saved state, then rename tab.
$ws.Activate()
$ws.Range("A7").Select()

$ws.Name = "N3mOpen"
